$wb = $excel.ActiveWorkbook

function Set-Txn($ws, $row, $a, $c, $d, $e, $f, $g, $h, $i, $j, $k, $n) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 3).Value = $c
    if ($d -ne $null) { $ws.Cells.Item($row, 4).Value = $d } else { $ws.Cells.Item($row, 4).Value = "" }
    if ($e -ne $null) { $ws.Cells.Item($row, 5).Value = $e } else { $ws.Cells.Item($row, 5).Value = "" }
    if ($f -ne $null) { $ws.Cells.Item($row, 6).Value = $f } else { $ws.Cells.Item($row, 6).Value = "" }
    if ($g -ne $null) { $ws.Cells.Item($row, 7).Value = $g } else { $ws.Cells.Item($row, 7).Value = "" }
    $ws.Cells.Item($row, 8).Value = $h
    if ($i -ne $null) { $ws.Cells.Item($row, 9).Value = $i } else { $ws.Cells.Item($row, 9).Value = "" }
    if ($j -ne $null) { $ws.Cells.Item($row, 10).Value = $j } else { $ws.Cells.Item($row, 10).Value = "" }
    if ($k -ne $null) { $ws.Cells.Item($row, 11).Value = $k } else { $ws.Cells.Item($row, 11).Value = "" }
    if ($n -ne $null) { $ws.Cells.Item($row, 14).Value = $n } else { $ws.Cells.Item($row, 14).Value = "" }
}

# "Memo blank" and "Also Memo blank" sheets share identical transaction data
foreach ($sheetName in @("Memo blank", "Also Memo blank")) {
    $ws = $wb.Worksheets.Item($sheetName)

    Set-Txn $ws 2 1 43102 'Wendy''s' 'Dining Out' $null $null -994.07 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 3 2 43102 'Walgreens' 'Pharmacy' 'General' $null -150.97 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 4 3 43102 'Amazon' 'Shopping' 'Online' $null -427.77 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 5 4 43104 'Walmart' 'Shopping' 'Household' $null -106.34 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 6 5 43104 'Amazon' 'Shopping' 'Online' $null -168.82 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 7 6 43105 'Neighborhood Market' 'Groceries' 'Farmer''s Market' $null -135.71 'Expense' 'Essential' $null 'Era C'
    Set-Txn $ws 8 7 43105 'Amazon' 'Shopping' 'Online' $null -898.86 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 9 8 43106 'CVS' 'Pharmacy' 'OTC' $null -729.9 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 10 9 43106 'Target' 'Shopping' 'Household' $null -336.88 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 11 10 43109 'Walgreens' 'Pharmacy' 'General' $null -608.97 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 12 11 43109 'Target' 'Shopping' 'Household' $null -759.35 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 13 12 43113 'Walgreens' 'Pharmacy' 'General' $null -478.88 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 14 13 43114 'Neighborhood Market' 'Groceries' 'Farmer''s Market' $null -159.92 'Expense' 'Essential' $null 'Era C'
    Set-Txn $ws 15 14 43114 'Target' 'Shopping' 'Household' $null -454.97 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 16 15 43115 'Amazon' 'Shopping' 'Online' $null -455.74 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 17 16 43116 'Amazon' 'Shopping' 'Online' $null -7.9 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 18 17 43119 'McDonald''s' 'Dining Out' $null $null -812.6 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 19 18 43119 'Neighborhood Market' 'Groceries' 'Farmer''s Market' $null -544.89 'Expense' 'Essential' $null 'Era C'
    Set-Txn $ws 20 19 43122 'Amazon' 'Shopping' 'Online' $null -68.91 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 21 20 43124 'McDonald''s' 'Dining Out' $null $null -907.43 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 22 21 43126 'Kroger' 'Groceries' 'Grocery Store' $null -828 'Expense' 'Essential' $null 'Era C'
    Set-Txn $ws 23 22 43126 'Walgreens' 'Pharmacy' 'General' $null -458.26 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 24 23 43126 'Amazon' 'Shopping' 'Online' $null -127.81 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 25 24 43127 'Amazon' 'Shopping' 'Online' $null -982.34 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 26 25 43128 'CVS' 'Pharmacy' 'OTC' $null -8.55 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 27 26 43130 'Wendy''s' 'Dining Out' $null $null -786.94 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 28 27 43130 'Walgreens' 'Pharmacy' 'General' $null -886.58 'Expense' 'Discretionary' $null 'Era C'
    Set-Txn $ws 29 28 43131 $null 'Transfer From' $null 'Checking' 12950.48 'Transfer' 'Transfer' $null 'Era C'

    $ws.Columns.Item(8).ColumnWidth = 10.166666666666666
    $ws.Columns.Item(11).ColumnWidth = 22.166666666666668
}

# "Memo not blank" and "Also Memo not blank" sheets share identical transaction data
foreach ($sheetName in @("Memo not blank", "Also Memo not blank")) {
    $ws = $wb.Worksheets.Item($sheetName)

    Set-Txn $ws 2 1 43109 'Target' 'Shopping' 'Household' $null 336.88 'Expense' 'Discretionary' 'Refund of 01/06/2018' 'Era C'
    Set-Txn $ws 3 2 43182 'Walmart' 'Shopping' 'Household' $null -484.39 'Expense' 'Discretionary' 'Memo 14' 'Era C'
    Set-Txn $ws 4 3 43243 'Target' 'Shopping' 'Household' $null 4.27 'Expense' 'Discretionary' 'Refund of 05/18/2018' 'Era C'
    Set-Txn $ws 5 4 43262 'McDonald''s' 'Dining Out' $null $null -727.42 'Expense' 'Discretionary' 'Memo 10' 'Era C'
    Set-Txn $ws 6 5 43262 'Walgreens' 'Pharmacy' 'General' $null -435.61 'Expense' 'Discretionary' 'Memo 18' 'Era C'

    $ws.Columns.Item(8).ColumnWidth = 10.166666666666666
    $ws.Columns.Item(11).ColumnWidth = 22.166666666666668
}
